$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Cover-letter date: "28 February 2018" -> "1 March 2018", with the
#    _GoBack bookmark relocated so it sits right after "1 March" (and
#    before " 2018"). The bookmark previously lived inside the title
#    paragraph (#15); it is removed from there naturally because that
#    whole span of text gets replaced in step (2) below.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1).Range
$p1.Find.Execute("28 February 2018", $true, $false, $false, $false, $false, $true, 1, $false, "1 March 2018", 2) | Out-Null

# Remove the old bookmark (currently inside paragraph 15's text) so we
# can re-create it at the new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-find paragraph 1 and place the bookmark right after "1 March".
$p1 = $d.Paragraphs.Item(1).Range
$p1.Find.Execute("1 March", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($p1.End, $p1.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ------------------------------------------------------------------
# 2) Title paragraph: drop ' as a Management Brief.' and add a period
#    before the closing curly quote; add "the" before "associate editor".
# ------------------------------------------------------------------
$p15 = $d.Paragraphs.Item(15).Range
$p15.Find.Execute(
    [char]0x201C + "Quantile Regression Estimates of Body Weight at Length in Walleye" + [char]0x201D + " as a Management Brief. I have incorporated the revisions suggested by you, associate editor,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    [char]0x201C + "Quantile Regression Estimates of Body Weight at Length in Walleye." + [char]0x201D + " I have incorporated the revisions suggested by you, the associate editor,",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) AE general comment #1.
# ------------------------------------------------------------------
$p49 = $d.Paragraphs.Item(49).Range
$p49.Find.Execute("incorporated into the methods of the manuscript.", $true, $false, $false, $false, $false, $true, 1, $false, "incorporated into the analysis of the data in this revision.", 2) | Out-Null

# ------------------------------------------------------------------
# 4) AE general comment #3: add sentence about software packages.
# ------------------------------------------------------------------
$p53 = $d.Paragraphs.Item(53).Range
$p53.Find.Execute(
    "in their analysis. Converting anyone",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in their analysis. I have also pointed out that quantile regression methods are available in commonly-used data-analysis software packages. Converting anyone",
    2) | Out-Null

# ------------------------------------------------------------------
# 5) "The reality, though, ..." paragraph: add ", scientists, ... and
#    reviewers" to the list of people.
# ------------------------------------------------------------------
$p55 = $d.Paragraphs.Item(55).Range
$p55.Find.Execute(
    "hopefully, moving forward, managers and authors will begin",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "hopefully, moving forward, managers, scientists, authors, and reviewers will begin",
    2) | Out-Null

# ------------------------------------------------------------------
# 6) Reviewer-1-general-comments wrap-up: add "with this revision".
# ------------------------------------------------------------------
$p69 = $d.Paragraphs.Item(69).Range
$p69.Find.Execute(
    "the code I include as a supplement",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the code I include with this revision as a supplement",
    2) | Out-Null

# ------------------------------------------------------------------
# 7) "Figure 1: Revised." paragraph: append explanatory sentence.
# ------------------------------------------------------------------
$p83 = $d.Paragraphs.Item(83).Range
$p83.Find.Execute(
    "Figure 1: Revised. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Figure 1: Revised. I have included a figure that has two panels, one for the reference data set and GA populations and one for reference dataset and SD populations. In greyscale, the confidence bands are difficult to differentiate, but evident. I have available the same image but in color that, if warranted, the AFS editorial board could use in place of the greyscale version included here. ",
    2) | Out-Null

# ------------------------------------------------------------------
# 8) "Line 190-191, 199-200" paragraph: "several years" -> "many years".
# ------------------------------------------------------------------
$p105 = $d.Paragraphs.Item(105).Range
$p105.Find.Execute("David Willis several years ago.", $true, $false, $false, $false, $false, $true, 1, $false, "David Willis many years ago.", 2) | Out-Null

Write-Output "done"
